$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Remove the "Adding Touch Designer to PATH" section entirely.
#    This covers six whole paragraphs:
#      - Heading3: "Adding Touch Designer to PATH"
#      - Normal:   "In order to run the Python program properly..."
#      - Heading4: "MacOS"
#      - Normal:   "Run: ..."
#      - Normal:   "Add the directory to the path: ..."
#      - Heading4: "Windows"
#    plus the text of the final paragraph ("Follow normal practices
#    for finding Touch Designer application. "), while keeping the
#    trailing single-space run that paragraph ended with -- that
#    space gets appended onto the prior ("...queue if many users
#    want to generate art. ") paragraph instead.
# -----------------------------------------------------------------

$anchorText = "Initially, it was thought that disabling real-time processing"

# Locate the paragraph index of the anchor ("...queue if many users
# want to generate art. ") paragraph -- the one immediately before
# the section we want to remove.
$targetIndex = 0
$idx = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$anchorText*") {
        $targetIndex = $idx
    }
    $idx = $idx + 1
}

if ($targetIndex -eq 0) {
    throw "Could not find paragraph index for anchor text"
}

$headingIndex = $targetIndex + 1

# Delete the six paragraphs that make up the "Adding Touch Designer to
# PATH" section (Heading3 through the Windows Heading4). Repeatedly
# delete the paragraph now sitting at $headingIndex, since later
# paragraphs' indices shift down automatically once one is removed.
for ($n = 0; $n -lt 6; $n++) {
    $p = $d.Paragraphs.Item($headingIndex)
    $p.Range.Delete()
}

# The paragraph now at $headingIndex is the former "Follow normal
# practices for finding Touch Designer application.  " paragraph.
# Remove it completely (text + paragraph mark) ...
$lastPara = $d.Paragraphs.Item($headingIndex)
$lastPara.Range.Delete()

# ... then re-append just the trailing single space run onto the end
# of the preceding ("queue if many users...") paragraph, matching
# what remained in the source paragraph after its heading text was
# removed.
$prevPara = $d.Paragraphs.Item($headingIndex - 1)
$prevPara.Range.InsertAfter(" ")

# -----------------------------------------------------------------
# 2. Update the cached "Last Revised" DATE field result in the
#    document header from "April 30, 2024" to "May 1, 2024".
# -----------------------------------------------------------------
$sec = $d.Sections.Item(1)
$header = $sec.Headers.Item(1)
$header.Range.Find.Execute("April 30, 2024", $true, $false, $false, $false, $false, $true, 1, $false, "May 1, 2024", 2) | Out-Null
